# Applies the cryptos list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving it as literal text
# (Price/Volume columns use plain-text numeric-looking strings such as
# "1.009" or "26.216.94"; without forcing a Text format Excel would
# silently convert them to floating point numbers).
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$updates = @(
    @{ Cell = 'D2'; Value = '26.216.94'; IsText = $false }
    @{ Cell = 'E2'; Value = '  -4.16%  '; IsText = $false }
    @{ Cell = 'D3'; Value = '1.658.99'; IsText = $false }
    @{ Cell = 'E3'; Value = '  -2.91%  '; IsText = $false }
    @{ Cell = 'D4'; Value = '1.009'; IsText = $false }
    @{ Cell = 'E4'; Value = '  +0.52%  '; IsText = $false }
    @{ Cell = 'D5'; Value = '218.13'; IsText = $false }
    @{ Cell = 'E5'; Value = '  -2.61%  '; IsText = $false }
    @{ Cell = 'D6'; Value = '0.5161'; IsText = $false }
    @{ Cell = 'E6'; Value = '  -3.37%  '; IsText = $false }
    @{ Cell = 'D7'; Value = '1.009'; IsText = $false }
    @{ Cell = 'E7'; Value = '  +0.54%  '; IsText = $false }
    @{ Cell = 'D8'; Value = '0.2571'; IsText = $false }
    @{ Cell = 'E8'; Value = '  -3.79%  '; IsText = $false }
    @{ Cell = 'D9'; Value = '0.06375'; IsText = $false }
    @{ Cell = 'E9'; Value = '  -3.53%  '; IsText = $false }
    @{ Cell = 'D10'; Value = '19.83'; IsText = $false }
    @{ Cell = 'E10'; Value = '  -5.36%  '; IsText = $false }
    @{ Cell = 'D11'; Value = '0.07778'; IsText = $false }
    @{ Cell = 'E11'; Value = '  +1.86%  '; IsText = $false }
    @{ Cell = 'D12'; Value = '1.670.07'; IsText = $false }
    @{ Cell = 'E12'; Value = '  -2.58%  '; IsText = $false }
    @{ Cell = 'B13'; Value = 'WrappedliquidstakedEther2.0'; IsText = $true }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; IsText = $true }
    @{ Cell = 'D13'; Value = '1.886.83'; IsText = $false }
    @{ Cell = 'E13'; Value = '  -3.01%  '; IsText = $false }
    @{ Cell = 'B14'; Value = 'Polkadot'; IsText = $true }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; IsText = $true }
    @{ Cell = 'D14'; Value = '4.294'; IsText = $false }
    @{ Cell = 'E14'; Value = '  -5.76%  '; IsText = $false }
    @{ Cell = 'D15'; Value = '0.5517'; IsText = $false }
    @{ Cell = 'E15'; Value = '  -4.32%  '; IsText = $false }
    @{ Cell = 'D16'; Value = '0.0₅8017'; IsText = $false }
    @{ Cell = 'E16'; Value = '  -1.93%  '; IsText = $false }
    @{ Cell = 'D17'; Value = '64.13'; IsText = $false }
    @{ Cell = 'E17'; Value = '  -5.45%  '; IsText = $false }
    @{ Cell = 'D18'; Value = '26.243.68'; IsText = $false }
    @{ Cell = 'E18'; Value = '  -4.14%  '; IsText = $false }
    @{ Cell = 'E19'; Value = '  +0.46%  '; IsText = $false }
    @{ Cell = 'D20'; Value = '209.87'; IsText = $false }
    @{ Cell = 'E20'; Value = '  -3.37%  '; IsText = $false }
    @{ Cell = 'D21'; Value = '4.378'; IsText = $false }
    @{ Cell = 'E21'; Value = '  -6.18%  '; IsText = $false }
    @{ Cell = 'D22'; Value = '10.05'; IsText = $false }
    @{ Cell = 'E22'; Value = '  -3.81%  '; IsText = $false }
    @{ Cell = 'D23'; Value = '5.874'; IsText = $false }
    @{ Cell = 'E23'; Value = '  -1.73%  '; IsText = $false }
    @{ Cell = 'D24'; Value = '1.009'; IsText = $false }
    @{ Cell = 'E24'; Value = '  +0.51%  '; IsText = $false }
    @{ Cell = 'D25'; Value = '143.61'; IsText = $false }
    @{ Cell = 'E25'; Value = '  +0.86%  '; IsText = $false }
    @{ Cell = 'E26'; Value = '  +2.55%  '; IsText = $false }
    @{ Cell = 'D27'; Value = '0.1160'; IsText = $false }
    @{ Cell = 'E27'; Value = '  -4.53%  '; IsText = $false }
    @{ Cell = 'D28'; Value = '6.953'; IsText = $false }
    @{ Cell = 'E28'; Value = '  -4.45%  '; IsText = $false }
    @{ Cell = 'D29'; Value = '15.70'; IsText = $false }
    @{ Cell = 'E29'; Value = '  -3.52%  '; IsText = $false }
    @{ Cell = 'D30'; Value = '0.05245'; IsText = $false }
    @{ Cell = 'E30'; Value = '  -2.91%  '; IsText = $false }
    @{ Cell = 'E31'; Value = '  -2.55%  '; IsText = $false }
    @{ Cell = 'D32'; Value = '3.360'; IsText = $false }
    @{ Cell = 'E32'; Value = '  -3.85%  '; IsText = $false }
    @{ Cell = 'E33'; Value = '  -6.33%  '; IsText = $false }
    @{ Cell = 'D34'; Value = '1.568'; IsText = $false }
    @{ Cell = 'E34'; Value = '  -4.42%  '; IsText = $false }
    @{ Cell = 'D35'; Value = '2.760'; IsText = $false }
    @{ Cell = 'E35'; Value = '  -4.07%  '; IsText = $false }
    @{ Cell = 'D36'; Value = '2.366'; IsText = $false }
    @{ Cell = 'E36'; Value = '  -1.93%  '; IsText = $false }
    @{ Cell = 'D37'; Value = '0.9235'; IsText = $false }
    @{ Cell = 'E37'; Value = '  -2.76%  '; IsText = $false }
    @{ Cell = 'D38'; Value = '0.5699'; IsText = $false }
    @{ Cell = 'E38'; Value = '  -2.75%  '; IsText = $false }
    @{ Cell = 'D39'; Value = '1.153.86'; IsText = $false }
    @{ Cell = 'E39'; Value = '  +10.48%  '; IsText = $false }
    @{ Cell = 'E40'; Value = '  -2.82%  '; IsText = $false }
    @{ Cell = 'D41'; Value = '1.009'; IsText = $false }
    @{ Cell = 'E41'; Value = '  +0.53%  '; IsText = $false }
    @{ Cell = 'D42'; Value = '0.8389'; IsText = $false }
    @{ Cell = 'E42'; Value = '  -0.32%  '; IsText = $false }
    @{ Cell = 'D43'; Value = '5.669'; IsText = $false }
    @{ Cell = 'E43'; Value = '  -3.31%  '; IsText = $false }
    @{ Cell = 'D44'; Value = '99.86'; IsText = $false }
    @{ Cell = 'E44'; Value = '  -0.94%  '; IsText = $false }
    @{ Cell = 'D45'; Value = '1.796.85'; IsText = $false }
    @{ Cell = 'E45'; Value = '  -3.06%  '; IsText = $false }
    @{ Cell = 'D46'; Value = '0.0₈110'; IsText = $false }
    @{ Cell = 'E46'; Value = '  +0.12%  '; IsText = $false }
    @{ Cell = 'D47'; Value = '0.4513'; IsText = $false }
    @{ Cell = 'E47'; Value = '  +0.03%  '; IsText = $false }
    @{ Cell = 'D48'; Value = '55.96'; IsText = $false }
    @{ Cell = 'E48'; Value = '  -3.45%  '; IsText = $false }
    @{ Cell = 'E49'; Value = '  +0.63%  '; IsText = $false }
    @{ Cell = 'D50'; Value = '7.889'; IsText = $false }
    @{ Cell = 'E50'; Value = '  -2.72%  '; IsText = $false }
    @{ Cell = 'D51'; Value = '0.05092'; IsText = $false }
    @{ Cell = 'E51'; Value = '  -2.86%  '; IsText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.IsText) {
        $range.Value = $u.Value
    } else {
        Set-TextCell $range $u.Value
    }
}
